$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53 (pushes existing rows 53.. down by one)
$ws.Rows("53:53").Insert()

# Insert a new row at position 76 (pushes existing rows 76.. down by one more)
$ws.Rows("76:76").Insert()

# Fill new row 53
$ws.Range("A53").Value = 3
$ws.Range("B53").Value = "Femacal de La Calera"
$ws.Range("C53").Value = "Coquimbo"
$ws.Range("D53").Value = 44679
$ws.Range("E53").Value = 5
$ws.Range("F53").Value = 100112026
$ws.Range("G53").Value = "Haba"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 38
$ws.Range("K53").Value = 22000
$ws.Range("L53").Value = 22000
$ws.Range("M53").Value = 22000
$ws.Range("N53").Value = "`$/malla 25 kilos"
$ws.Range("O53").Value = "Provincia de Limarí"
$ws.Range("P53").Value = 880
$ws.Range("Q53").Value = 25
$ws.Range("R53").Value = "Hortaliza"

# Fill new row 76
$ws.Range("A76").Value = 3
$ws.Range("B76").Value = "Femacal de La Calera"
$ws.Range("C76").Value = "Coquimbo"
$ws.Range("D76").Value = 44680
$ws.Range("E76").Value = 5
$ws.Range("F76").Value = 100112026
$ws.Range("G76").Value = "Haba"
$ws.Range("H76").Value = "Sin especificar"
$ws.Range("I76").Value = "Primera"
$ws.Range("J76").Value = 38
$ws.Range("K76").Value = 21000
$ws.Range("L76").Value = 21000
$ws.Range("M76").Value = 21000
$ws.Range("N76").Value = "`$/malla 25 kilos"
$ws.Range("O76").Value = "Provincia de Limarí"
$ws.Range("P76").Value = 840
$ws.Range("Q76").Value = 25
$ws.Range("R76").Value = "Hortaliza"

